{"js": "// Office.js (Word JavaScript API) script implementing the commit:\n// \"Added NCI SVN location for test data dump.\"\n//\n// The underlying change, as shown by the canonical OOXML diff, does two\n// kinds of things:\n//   1) Collapses a handful of \"spell-checked\" run splits (runs that were\n//      broken apart around a <w:proofErr> pair) back into a single run,\n//      with the proofErr markers removed entirely.\n//   2) Replaces the \"Import dump located at ...\" prerequisite paragraph\n//      with three new paragraphs that add the Oracle/MySQL SVN dump\n//      locations.\n//\n// Because Office.js has no direct \"delete this element\" primitive for\n// <w:proofErr/>, every edit below is performed by replacing the OOXML of\n// the whole affected paragraph with freshly authored OOXML that has no\n// proofErr runs left in it \u2014 this guarantees a clean result structurally\n// equivalent to the target diff.\n\nconst PKG_OPEN =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>';\nconst PKG_CLOSE =\n  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\n// Replace a paragraph (identified by its Range) with one or more raw\n// <w:p>...</w:p> paragraph fragments.\nfunction replaceParagraphWithOoxml(range, paragraphsXml) {\n  range.insertOoxml(PKG_OPEN + paragraphsXml + PKG_CLOSE, Word.InsertLocation.replace);\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 1) \"... with short title REPORT_LOADING_Reports_get_loaded_successfully\"\n//    merge the two runs that had been split by a proofErr pair.\n// ---------------------------------------------------------------------\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"with short title REPORT_LOADING_Reports_get_loaded_successfully\") !== -1) {\n    const xml =\n      '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr>' +\n      '<w:spacing w:after=\"0\"/></w:pPr>' +\n      '<w:r><w:t xml:space=\"preserve\">Select Test case ID </w:t></w:r>' +\n      '<w:r><w:t>9590</w:t></w:r>' +\n      '<w:r><w:t xml:space=\"preserve\"> with short title REPORT_LOADING_Reports_get_loaded_successfully</w:t></w:r>' +\n      '</w:p>';\n    replaceParagraphWithOoxml(paragraphs.items[i].getRange(), xml);\n    break;\n  }\n}\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2) \"Import dump located at ...\" -> three paragraphs with the Oracle /\n//    MySQL NCI SVN dump locations.\n// ---------------------------------------------------------------------\nparagraphs.load(\"items/text\");\nawait context.sync();\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Import dump located at\") !== -1) {\n    const xml =\n      '<w:p><w:pPr><w:spacing w:after=\"0\"/></w:pPr>' +\n      '<w:r><w:t xml:space=\"preserve\">Import latest dump located at </w:t></w:r></w:p>' +\n      '<w:p><w:pPr><w:spacing w:after=\"0\"/></w:pPr>' +\n      '<w:r><w:t>Oracle: https://ncisvn.nci.nih.gov/svn/catissue_persistent/caTissue Database Dump/v2.0/Oracle</w:t></w:r></w:p>' +\n      '<w:p><w:pPr><w:spacing w:after=\"0\"/></w:pPr>' +\n      '<w:r><w:t xml:space=\"preserve\">MySQL: https://ncisvn.nci.nih.gov/svn/catissue_persistent/caTissue Database Dump/v2.0/MySQL and deploy application </w:t></w:r>' +\n      '<w:r><w:t>and deploy application with caTIES configuration.</w:t></w:r></w:p>';\n    replaceParagraphWithOoxml(paragraphs.items[i].getRange(), xml);\n    break;\n  }\n}\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 3) \"ant run_report_loader_server\" -> single bold run.\n// ---------------------------------------------------------------------\nparagraphs.load(\"items/text\");\nawait context.sync();\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"run_report_loader_server\") !== -1) {\n    const xml =\n      '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>' +\n      '<w:spacing w:after=\"0\"/></w:pPr>' +\n      '<w:r><w:t>Execute \\u201c</w:t></w:r>' +\n      '<w:r><w:rPr><w:b/></w:rPr><w:t>ant run_report_loader_server</w:t></w:r>' +\n      '<w:r><w:t>\\u201d target. (Refer Expected Output).</w:t></w:r>' +\n      '</w:p>';\n    replaceParagraphWithOoxml(paragraphs.items[i].getRange(), xml);\n    break;\n  }\n}\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 4) \"ant stop_report_loader_server\" -> single bold run.\n// ---------------------------------------------------------------------\nparagraphs.load(\"items/text\");\nawait context.sync();\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"stop_report_loader_server\") !== -1) {\n    const xml =\n      '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>' +\n      '<w:spacing w:after=\"0\"/></w:pPr>' +\n      '<w:r><w:t>Execute \\u201c</w:t></w:r>' +\n      '<w:r><w:rPr><w:b/></w:rPr><w:t>ant stop_report_loader_server</w:t></w:r>' +\n      '<w:r><w:t>\\u201d on a different command prompt from REPORT_LOADER_HOME directory. (Refer Expected Output)</w:t></w:r>' +\n      '</w:p>';\n    replaceParagraphWithOoxml(paragraphs.items[i].getRange(), xml);\n    break;\n  }\n}\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 5) \"S.No.\" table header cell appears 3 times -> merge into single bold run.\n// ---------------------------------------------------------------------\nparagraphs.load(\"items/text\");\nawait context.sync();\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"S.No.\") {\n    const xml =\n      '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr>' +\n      '<w:r><w:rPr><w:b/></w:rPr><w:t>S.No.</w:t></w:r>' +\n      '</w:p>';\n    replaceParagraphWithOoxml(paragraphs.items[i].getRange(), xml);\n  }\n}\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 6) \"Last Name: Synoptest (No match)\" -> merge into single run.\n// ---------------------------------------------------------------------\nparagraphs.load(\"items/text\");\nawait context.sync();\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Last Name: Synoptest (No match)\") {\n    const xml = '<w:p><w:r><w:t>Last Name: Synoptest (No match)</w:t></w:r></w:p>';\n    replaceParagraphWithOoxml(paragraphs.items[i].getRange(), xml);\n    break;\n  }\n}\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 7) \"5) Report loader server stops. ...\" -> merge the deploycaties.properties\n//    runs (previously split by two proofErr pairs) into a single run.\n// ---------------------------------------------------------------------\nparagraphs.load(\"items/text\");\nawait context.sync();\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"5) Report loader server stops.\") !== -1) {\n    const xml =\n      '<w:p><w:pPr><w:spacing w:after=\"0\"/></w:pPr>' +\n      '<w:r><w:t xml:space=\"preserve\">5) Report loader server stops. The reports will be visible to super administrator, ' +\n      'administrator of the site specified in deploycaties.properties file, PI or PC of the collection protocol ' +\n      'specified in the deploycaties.properties file. </w:t></w:r>' +\n      '</w:p>';\n    replaceParagraphWithOoxml(paragraphs.items[i].getRange(), xml);\n    break;\n  }\n}\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script implementing the commit:\n# \"Added NCI SVN location for test data dump.\"\n#\n# Mirrors edit.js: every affected paragraph's Range is replaced, in one\n# shot, via InsertXML with a freshly authored <w:p> (or sequence of\n# <w:p> elements). That is the cleanest way to both (a) merge runs that\n# had been split apart by a <w:proofErr/> pair \u2014 with the proofErr\n# markers dropped entirely \u2014 and (b) turn the \"Import dump located at\n# ...\" paragraph into three paragraphs carrying the Oracle/MySQL SVN\n# dump locations.\n\n$d = $word.ActiveDocument\n$W_NS = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n\nfunction Get-CleanText($range) {\n  # Paragraph.Range.Text includes the trailing paragraph mark (CR, 0x0D)\n  # and, for the last paragraph in a table cell, a cell-mark (0x07).\n  # Strip both so plain string comparisons work.\n  return $range.Text.TrimEnd([char]13, [char]7)\n}\n\n# ---------------------------------------------------------------------\n# 1) \"... with short title REPORT_LOADING_Reports_get_loaded_successfully\"\n#    merge the two runs that had been split by a proofErr pair.\n# ---------------------------------------------------------------------\nforeach ($p in $d.Paragraphs) {\n  if ($p.Range.Text -like \"*with short title REPORT_LOADING_Reports_get_loaded_successfully*\") {\n    $xml = \"<w:p $W_NS><w:pPr><w:pStyle w:val=`\"ListParagraph`\"/><w:numPr><w:ilvl w:val=`\"0`\"/><w:numId w:val=`\"2`\"/></w:numPr>\" +\n           \"<w:spacing w:after=`\"0`\"/></w:pPr>\" +\n           \"<w:r><w:t xml:space=`\"preserve`\">Select Test case ID </w:t></w:r>\" +\n           \"<w:r><w:t>9590</w:t></w:r>\" +\n           \"<w:r><w:t xml:space=`\"preserve`\"> with short title REPORT_LOADING_Reports_get_loaded_successfully</w:t></w:r>\" +\n           \"</w:p>\"\n    $p.Range.InsertXML($xml)\n    break\n  }\n}\n\n# ---------------------------------------------------------------------\n# 2) \"Import dump located at ...\" -> three paragraphs with the Oracle /\n#    MySQL NCI SVN dump locations.\n# ---------------------------------------------------------------------\nforeach ($p in $d.Paragraphs) {\n  if ($p.Range.Text -like \"*Import dump located at*\") {\n    $xml = \"<w:p $W_NS><w:pPr><w:spacing w:after=`\"0`\"/></w:pPr>\" +\n           \"<w:r><w:t xml:space=`\"preserve`\">Import latest dump located at </w:t></w:r></w:p>\" +\n           \"<w:p $W_NS><w:pPr><w:spacing w:after=`\"0`\"/></w:pPr>\" +\n           \"<w:r><w:t>Oracle: https://ncisvn.nci.nih.gov/svn/catissue_persistent/caTissue Database Dump/v2.0/Oracle</w:t></w:r></w:p>\" +\n           \"<w:p $W_NS><w:pPr><w:spacing w:after=`\"0`\"/></w:pPr>\" +\n           \"<w:r><w:t xml:space=`\"preserve`\">MySQL: https://ncisvn.nci.nih.gov/svn/catissue_persistent/caTissue Database Dump/v2.0/MySQL and deploy application </w:t></w:r>\" +\n           \"<w:r><w:t>and deploy application with caTIES configuration.</w:t></w:r></w:p>\"\n    $p.Range.InsertXML($xml)\n    break\n  }\n}\n\n# ---------------------------------------------------------------------\n# 3) \"ant run_report_loader_server\" -> single bold run.\n# ---------------------------------------------------------------------\nforeach ($p in $d.Paragraphs) {\n  if ($p.Range.Text -like \"*run_report_loader_server*\") {\n    $xml = \"<w:p $W_NS><w:pPr><w:pStyle w:val=`\"ListParagraph`\"/><w:numPr><w:ilvl w:val=`\"0`\"/><w:numId w:val=`\"1`\"/></w:numPr>\" +\n           \"<w:spacing w:after=`\"0`\"/></w:pPr>\" +\n           \"<w:r><w:t>Execute \" + [char]8220 + \"</w:t></w:r>\" +\n           \"<w:r><w:rPr><w:b/></w:rPr><w:t>ant run_report_loader_server</w:t></w:r>\" +\n           \"<w:r><w:t>\" + [char]8221 + \" target. (Refer Expected Output).</w:t></w:r>\" +\n           \"</w:p>\"\n    $p.Range.InsertXML($xml)\n    break\n  }\n}\n\n# ---------------------------------------------------------------------\n# 4) \"ant stop_report_loader_server\" -> single bold run.\n# ---------------------------------------------------------------------\nforeach ($p in $d.Paragraphs) {\n  if ($p.Range.Text -like \"*stop_report_loader_server*\") {\n    $xml = \"<w:p $W_NS><w:pPr><w:pStyle w:val=`\"ListParagraph`\"/><w:numPr><w:ilvl w:val=`\"0`\"/><w:numId w:val=`\"1`\"/></w:numPr>\" +\n           \"<w:spacing w:after=`\"0`\"/></w:pPr>\" +\n           \"<w:r><w:t>Execute \" + [char]8220 + \"</w:t></w:r>\" +\n           \"<w:r><w:rPr><w:b/></w:rPr><w:t>ant stop_report_loader_server</w:t></w:r>\" +\n           \"<w:r><w:t>\" + [char]8221 + \" on a different command prompt from REPORT_LOADER_HOME directory. (Refer Expected Output)</w:t></w:r>\" +\n           \"</w:p>\"\n    $p.Range.InsertXML($xml)\n    break\n  }\n}\n\n# ---------------------------------------------------------------------\n# 5) \"S.No.\" table header cell appears 3 times -> merge into single bold run.\n# ---------------------------------------------------------------------\nforeach ($p in $d.Paragraphs) {\n  if ((Get-CleanText $p.Range) -eq \"S.No.\") {\n    $xml = \"<w:p $W_NS><w:pPr><w:rPr><w:b/></w:rPr></w:pPr>\" +\n           \"<w:r><w:rPr><w:b/></w:rPr><w:t>S.No.</w:t></w:r>\" +\n           \"</w:p>\"\n    $p.Range.InsertXML($xml)\n  }\n}\n\n# ---------------------------------------------------------------------\n# 6) \"Last Name: Synoptest (No match)\" -> merge into single run.\n# ---------------------------------------------------------------------\nforeach ($p in $d.Paragraphs) {\n  if ((Get-CleanText $p.Range) -eq \"Last Name: Synoptest (No match)\") {\n    $xml = \"<w:p $W_NS><w:r><w:t>Last Name: Synoptest (No match)</w:t></w:r></w:p>\"\n    $p.Range.InsertXML($xml)\n    break\n  }\n}\n\n# ---------------------------------------------------------------------\n# 7) \"5) Report loader server stops. ...\" -> merge the deploycaties.properties\n#    runs (previously split by two proofErr pairs) into a single run.\n# ---------------------------------------------------------------------\nforeach ($p in $d.Paragraphs) {\n  if ($p.Range.Text -like \"*5) Report loader server stops.*\") {\n    $xml = \"<w:p $W_NS><w:pPr><w:spacing w:after=`\"0`\"/></w:pPr>\" +\n           \"<w:r><w:t xml:space=`\"preserve`\">5) Report loader server stops. The reports will be visible to super administrator, \" +\n           \"administrator of the site specified in deploycaties.properties file, PI or PC of the collection protocol \" +\n           \"specified in the deploycaties.properties file. </w:t></w:r>\" +\n           \"</w:p>\"\n    $p.Range.InsertXML($xml)\n    break\n  }\n}\n"}
